$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.081.27"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "2.105.85"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").Value = "'345.06"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("D8").Value = "'0.4456"
$ws.Range("E8").Value = "  -1.00%  "
$ws.Range("D9").Value = "'0.09512"
$ws.Range("E9").Value = "  +4.75%  "
$ws.Range("D10").Value = "'52.50"
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "'25.38"
$ws.Range("E12").Value = "  +4.18%  "
$ws.Range("D13").Value = "2.106.09"
$ws.Range("D14").Value = "'6.751"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "'8.117"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "'99.71"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "'0.00001171"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "'20.73"
$ws.Range("E19").Value = "  +7.17%  "
$ws.Range("D20").Value = "'0.06702"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'1.006"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "'6.207"
$ws.Range("E22").Value = "  -3.15%  "
$ws.Range("D23").Value = "30.166.15"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("D24").Value = "'12.74"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "'2.318"
$ws.Range("E25").Value = "  -2.41%  "
$ws.Range("D26").Value = "2.353.68"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "'22.08"
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").Value = "'164.44"
$ws.Range("E28").Value = "  -0.27%  "
$ws.Range("D29").Value = "'2.544"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'133.74"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").Value = "'1.162"
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("D33").Value = "'1.633"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'6.261"
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").Value = "'3.942"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'6.183"
$ws.Range("E36").Value = "  +4.66%  "
$ws.Range("D37").Value = "'10.17"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").Value = "'0.02581"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").Value = "'0.06790"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").Value = "'0.2289"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").Value = "'0.6967"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").Value = "'12.53"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "'1.307"
$ws.Range("E43").Value = "  +3.72%  "
$ws.Range("D44").Value = "'0.6704"
$ws.Range("E44").Value = "  +4.29%  "
$ws.Range("D45").Value = "'14.30"
$ws.Range("E45").Value = "  -5.21%  "
$ws.Range("D46").Value = "'2.291"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").Value = "'0.00000000357"
$ws.Range("E48").Value = "  -4.00%  "
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").Value = "'82.24"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("E51").Value = "  -1.83%  "
